$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

$row1 = $tbl.ListRows.Add()
$row1.Range.Cells.Item(1, 1).Value = "LookupValue"
$row1.Range.Cells.Item(1, 2).Value = "Permissions"

$row2 = $tbl.ListRows.Add()
$row2.Range.Cells.Item(1, 1).Value = "LookupValue"
$row2.Range.Cells.Item(1, 2).Value = "Copy Document"

$ws.Range("C16").Select()
